{"js": "// The practice sheet keeps the date in the first paragraph and 100 answer\n// cells (20 rows x 5 columns, in reading order) inside the single table\n// that follows it. Refresh the date and regenerate every problem/answer\n// pair while leaving all paragraph/run formatting untouched.\n\nconst body = context.document.body;\n\n// --- Update the date line -------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2025-04-05 Saturday\", \"Replace\");\n\n// --- Update every answer cell in the table --------------------------------\nconst values = [\n  \"85+3=88\",\n  \"57+7=64\",\n  \"49-30=19\",\n  \"63-53=10\",\n  \"90-5=85\",\n  \"56+11=67\",\n  \"24+58=82\",\n  \"76-30=46\",\n  \"3+46=49\",\n  \"29+18=47\",\n  \"85+12=97\",\n  \"68+29=97\",\n  \"83-71=12\",\n  \"75+13=88\",\n  \"71+28=99\",\n  \"56+40=96\",\n  \"62+14=76\",\n  \"12+72=84\",\n  \"84-34=50\",\n  \"25+61=86\",\n  \"71-61=10\",\n  \"54+3=57\",\n  \"96-10=86\",\n  \"80-53=27\",\n  \"61+35=96\",\n  \"15+20=35\",\n  \"27+67=94\",\n  \"91-55=36\",\n  \"64+30=94\",\n  \"3+84=87\",\n  \"53+18=71\",\n  \"2+55=57\",\n  \"82-46=36\",\n  \"31+3=34\",\n  \"75-40=35\",\n  \"39-35=4\",\n  \"66-14=52\",\n  \"46-16=30\",\n  \"74-6=68\",\n  \"72-71=1\",\n  \"52+23=75\",\n  \"64+29=93\",\n  \"92+2=94\",\n  \"39+42=81\",\n  \"5+33=38\",\n  \"17-9=8\",\n  \"87-65=22\",\n  \"22-6=16\",\n  \"8+50=58\",\n  \"81-64=17\",\n  \"5+39=44\",\n  \"89-3=86\",\n  \"81+1=82\",\n  \"49+48=97\",\n  \"17+42=59\",\n  \"91-66=25\",\n  \"54-54=0\",\n  \"12+43=55\",\n  \"94-48=46\",\n  \"67-11=56\",\n  \"84-76=8\",\n  \"25+17=42\",\n  \"96-9=87\",\n  \"4+29=33\",\n  \"40-8=32\",\n  \"37+24=61\",\n  \"23-15=8\",\n  \"7+25=32\",\n  \"28-25=3\",\n  \"39-31=8\",\n  \"99-97=2\",\n  \"83-26=57\",\n  \"20+42=62\",\n  \"41-6=35\",\n  \"70-15=55\",\n  \"81-14=67\",\n  \"29+48=77\",\n  \"49+17=66\",\n  \"49-28=21\",\n  \"8+69=77\",\n  \"16+21=37\",\n  \"27+17=44\",\n  \"82-9=73\",\n  \"15+63=78\",\n  \"7-4=3\",\n  \"0+47=47\",\n  \"83-49=34\",\n  \"39-14=25\",\n  \"22+38=60\",\n  \"72-72=0\",\n  \"27+30=57\",\n  \"48-6=42\",\n  \"73+5=78\",\n  \"21+31=52\",\n  \"42+2=44\",\n  \"59-16=43\",\n  \"26+52=78\",\n  \"43+24=67\",\n  \"82-62=20\",\n  \"31+31=62\"\n];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = 20;\nconst cols = 5;\n\n// Load each cell's first paragraph so we can grab its range afterwards.\nconst cellParagraphs = [];\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const cellParas = table.getCell(r, c).body.paragraphs;\n    cellParas.load(\"items\");\n    cellParagraphs.push(cellParas);\n  }\n}\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const p = cellParagraphs[idx].items[0];\n    p.getRange().insertText(values[idx], \"Replace\");\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line at the top of the document.\n$d = $word.ActiveDocument\n$d.Paragraphs(1).Range.Text = \"2025-04-05 Saturday\"\n\n# Update every answer cell in the practice table (5 columns x 20 rows, in\n# reading order) to the new set of addition/subtraction problems.\n$values = @(\n    \"85+3=88\",\n    \"57+7=64\",\n    \"49-30=19\",\n    \"63-53=10\",\n    \"90-5=85\",\n    \"56+11=67\",\n    \"24+58=82\",\n    \"76-30=46\",\n    \"3+46=49\",\n    \"29+18=47\",\n    \"85+12=97\",\n    \"68+29=97\",\n    \"83-71=12\",\n    \"75+13=88\",\n    \"71+28=99\",\n    \"56+40=96\",\n    \"62+14=76\",\n    \"12+72=84\",\n    \"84-34=50\",\n    \"25+61=86\",\n    \"71-61=10\",\n    \"54+3=57\",\n    \"96-10=86\",\n    \"80-53=27\",\n    \"61+35=96\",\n    \"15+20=35\",\n    \"27+67=94\",\n    \"91-55=36\",\n    \"64+30=94\",\n    \"3+84=87\",\n    \"53+18=71\",\n    \"2+55=57\",\n    \"82-46=36\",\n    \"31+3=34\",\n    \"75-40=35\",\n    \"39-35=4\",\n    \"66-14=52\",\n    \"46-16=30\",\n    \"74-6=68\",\n    \"72-71=1\",\n    \"52+23=75\",\n    \"64+29=93\",\n    \"92+2=94\",\n    \"39+42=81\",\n    \"5+33=38\",\n    \"17-9=8\",\n    \"87-65=22\",\n    \"22-6=16\",\n    \"8+50=58\",\n    \"81-64=17\",\n    \"5+39=44\",\n    \"89-3=86\",\n    \"81+1=82\",\n    \"49+48=97\",\n    \"17+42=59\",\n    \"91-66=25\",\n    \"54-54=0\",\n    \"12+43=55\",\n    \"94-48=46\",\n    \"67-11=56\",\n    \"84-76=8\",\n    \"25+17=42\",\n    \"96-9=87\",\n    \"4+29=33\",\n    \"40-8=32\",\n    \"37+24=61\",\n    \"23-15=8\",\n    \"7+25=32\",\n    \"28-25=3\",\n    \"39-31=8\",\n    \"99-97=2\",\n    \"83-26=57\",\n    \"20+42=62\",\n    \"41-6=35\",\n    \"70-15=55\",\n    \"81-14=67\",\n    \"29+48=77\",\n    \"49+17=66\",\n    \"49-28=21\",\n    \"8+69=77\",\n    \"16+21=37\",\n    \"27+17=44\",\n    \"82-9=73\",\n    \"15+63=78\",\n    \"7-4=3\",\n    \"0+47=47\",\n    \"83-49=34\",\n    \"39-14=25\",\n    \"22+38=60\",\n    \"72-72=0\",\n    \"27+30=57\",\n    \"48-6=42\",\n    \"73+5=78\",\n    \"21+31=52\",\n    \"42+2=44\",\n    \"59-16=43\",\n    \"26+52=78\",\n    \"43+24=67\",\n    \"82-62=20\",\n    \"31+31=62\"\n)\n\n$t = $d.Tables(1)\n$cols = $t.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$idx]\n        $idx++\n    }\n}\n"}
